# Atualização de bases das ligas, do dia: 14-05-2024 às 20:19
#
# This script reproduces, via the Excel object model, the same net effect
# the authoritative OOXML diff shows for "Mexico Liga de Expansion.xlsx":
#   1) Rows 91 and 92 (match ids 6924569 / 6924568) have their full data
#      (every column except the running "id" in column A) swapped.
#   2) Rows 186 and 187 (match ids 7648958 / 7648957) likewise have their
#      full data swapped.
#   3) A brand new match record (id 8185481) is appended as row 248,
#      extending the sheet dimension from A1:AB247 to A1:AB248.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param([int]$Row1, [int]$Row2)
    # Columns B..AB (2..28) hold the data that should be exchanged between
    # the two rows; column A (the sequential id) stays put.
    # NOTE: positional (not named) arguments are used when calling this
    # function, since named parameter binding is not reliable here.
    for ($col = 2; $col -le 28; $col++) {
        $cell1 = $ws.Cells.Item($Row1, $col)
        $cell2 = $ws.Cells.Item($Row2, $col)
        $tmp = $cell1.Value()
        $cell1.Value = $cell2.Value()
        $cell2.Value = $tmp
    }
}

# 1) Swap rows 91 / 92
Swap-RowData 91 92

# 2) Swap rows 186 / 187
Swap-RowData 186 187

# 3) Append the new row 248, copying the A/D column formatting (bold
#    centered id style, and the yyyy-mm-dd date format) from the last
#    existing data row (247) so the new row matches the sheet's styling.
$ws.Cells.Item(247, 1).Copy()
$ws.Cells.Item(248, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(247, 4).Copy()
$ws.Cells.Item(248, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(248, 1).Value = 246
$ws.Cells.Item(248, 2).Value = 8185481
$ws.Cells.Item(248, 3).Value = "Mexico Liga de Expansion"
$ws.Cells.Item(248, 4).Value = 45424.83333333334
$ws.Cells.Item(248, 5).Value = "Universidad Guadalajara"
$ws.Cells.Item(248, 6).Value = "Atlante"
$ws.Cells.Item(248, 7).Value = 1
$ws.Cells.Item(248, 8).Value = 2
$ws.Cells.Item(248, 9).Value = "A"
$ws.Cells.Item(248, 10).Value = 2.625
$ws.Cells.Item(248, 11).Value = 3
$ws.Cells.Item(248, 12).Value = 2.6
$ws.Cells.Item(248, 13).Value = 2.45
$ws.Cells.Item(248, 14).Value = 3
$ws.Cells.Item(248, 15).Value = 3.1
$ws.Cells.Item(248, 16).Value = -0.25
$ws.Cells.Item(248, 17).Value = 2.025
$ws.Cells.Item(248, 18).Value = 1.775
$ws.Cells.Item(248, 19).Value = 2
$ws.Cells.Item(248, 20).Value = 2
$ws.Cells.Item(248, 21).Value = 1.8
$ws.Cells.Item(248, 22).Value = -1
$ws.Cells.Item(248, 23).Value = -1
$ws.Cells.Item(248, 24).Value = 2.1
$ws.Cells.Item(248, 25).Value = -1
$ws.Cells.Item(248, 26).Value = 0.7749999999999999
$ws.Cells.Item(248, 27).Value = 1
$ws.Cells.Item(248, 28).Value = -1

Write-Output "Edit applied: swapped rows 91/92, swapped rows 186/187, appended row 248."
